# Update the "想去人数" (interested-count) values in column F
# across the "展览", "演出" and "全部类型" worksheets,
# matching the canonical OOXML diff (output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 19
$wsExpo.Range("F6").Value = 83
$wsExpo.Range("F7").Value = 890
$wsExpo.Range("F8").Value = 56
$wsExpo.Range("F9").Value = 7042
$wsExpo.Range("F11").Value = 153
$wsExpo.Range("F12").Value = 6526
$wsExpo.Range("F13").Value = 132
$wsExpo.Range("F15").Value = 4471
$wsExpo.Range("F18").Value = 4514
$wsExpo.Range("F19").Value = 12
$wsExpo.Range("F21").Value = 255
$wsExpo.Range("F30").Value = 8112
$wsExpo.Range("F32").Value = 1398
$wsExpo.Range("F34").Value = 708
$wsExpo.Range("F37").Value = 987
$wsExpo.Range("F39").Value = 1646
$wsExpo.Range("F40").Value = 213
$wsExpo.Range("F41").Value = 953
$wsExpo.Range("F43").Value = 4169
$wsExpo.Range("F45").Value = 28
$wsExpo.Range("F46").Value = 115
$wsExpo.Range("F48").Value = 839
$wsExpo.Range("F49").Value = 1116
$wsExpo.Range("F50").Value = 18

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F19").Value = 874

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 19
$wsAll.Range("F9").Value = 83
$wsAll.Range("F11").Value = 890
$wsAll.Range("F12").Value = 56
$wsAll.Range("F13").Value = 7042
$wsAll.Range("F15").Value = 153
$wsAll.Range("F16").Value = 6526
$wsAll.Range("F17").Value = 132
$wsAll.Range("F19").Value = 4471
$wsAll.Range("F22").Value = 4514
$wsAll.Range("F24").Value = 255
$wsAll.Range("F31").Value = 8112
$wsAll.Range("F33").Value = 1398
$wsAll.Range("F35").Value = 708
$wsAll.Range("F38").Value = 987
$wsAll.Range("F39").Value = 1646
$wsAll.Range("F40").Value = 213
$wsAll.Range("F41").Value = 953
$wsAll.Range("F43").Value = 4169
$wsAll.Range("F45").Value = 28
$wsAll.Range("F46").Value = 115
$wsAll.Range("F48").Value = 839
$wsAll.Range("F49").Value = 1116
$wsAll.Range("F50").Value = 18

